$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row (row 12) of trade data, matching existing columns A-H
$ws.Range("A12").Value = 9219.31
$ws.Range("B12").Value = 9138.89
$ws.Range("C12").Value = 105.78
$ws.Range("D12").Value = 106.71
$ws.Range("E12").Value = $false
$ws.Range("F12").Value = 0.88
$ws.Range("G12").Value = 42620.766192129631
$ws.Range("H12").Value = $true

# Column G uses a date/time display format, same as the rows above it
$ws.Range("G12").NumberFormat = "m/d/yy h:mm"
